$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the "b.md" entry, now ready for handoff ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-13 18:45:11"

# --- zh-cn sheet: row 3 is the "b.md" entry ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-13 18:44:59"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/7fc3806ca1a0cf1e73455bd95cdc61beb041d299/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/2e5432184113f3f45db89f443867f7858c01f9d1/e2e/b.md."
$zhcn.Range("P1").ColumnWidth = 39.15

# --- de-de sheet: row 3 is the "b.md" entry ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
# Content Duplicate flips from True to False; copy the text-typed "False" cell
# from the row above so the written value keeps its text data type instead of
# being auto-coerced into a native boolean.
$dede.Range("F2").Copy($dede.Range("F3"))
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-13 18:45:11"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/7fc3806ca1a0cf1e73455bd95cdc61beb041d299/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/2e5432184113f3f45db89f443867f7858c01f9d1/e2e/b.md."
$dede.Range("P1").ColumnWidth = 39.15
